$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")
$ws.Range("A1").Value = "Data#"
$ws.Range("B1").Value = "Test_Description"
$ws.Range("C1").Value = "user"
$ws.Range("D1").Value = "selectType"
$ws.Range("E1").Value = "search"
$ws.Range("F1").Value = "rationale"
$ws.Range("G1").Value = "observation"
$ws.Range("H1").Value = "enterText"
$ws.Range("I1").Value = "trait"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Check Delete Button disabled in tabbed view GF (GF with Evidence cannot be deleted)"
$ws.Range("C2").Value = "Pillai, Nisha"
$ws.Range("D2").Value = "Genetic Feature"
$ws.Range("E2").Value = "Bharitkar S, Mendel"
$ws.Range("F2").Value = "test"
$ws.Range("G2").Value = "test observation"
$ws.Range("H2").Value = "ath-MIR156a"
$ws.Range("I2").Value = "biomass yield [en;XX;1]"
